# Updated cryptos list on Mon May  1 07:36:37 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '28.674.23'
$ws.Range('E2').Value = '  -3.07%  '
$ws.Range('D3').Value = '1.851.76'
$ws.Range('E3').Value = '  -3.69%  '
$ws.Range('E4').Value = '  -1.04%  '
$ws.Range('D5').Value = '''334.45'
$ws.Range('E5').Value = '  +2.70%  '
$ws.Range('D6').Value = '''1.002'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('D7').Value = '''0.4659'
$ws.Range('E7').Value = '  -3.22%  '
$ws.Range('D8').Value = '''0.3914'
$ws.Range('E8').Value = '  -3.55%  '
$ws.Range('D9').Value = '''46.37'
$ws.Range('E9').Value = '  -2.68%  '
$ws.Range('D10').Value = '''0.07916'
$ws.Range('E10').Value = '  -3.66%  '
$ws.Range('D11').Value = '''0.9853'
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('D12').Value = '''22.36'
$ws.Range('E12').Value = '  -5.46%  '
$ws.Range('D13').Value = '1.840.92'
$ws.Range('E13').Value = '  -4.61%  '
$ws.Range('D14').Value = '''5.849'
$ws.Range('E14').Value = '  -3.68%  '
$ws.Range('D15').Value = '''6.998'
$ws.Range('E15').Value = '  -3.64%  '
$ws.Range('D16').Value = '''0.06850'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '''87.72'
$ws.Range('E17').Value = '  -4.25%  '
$ws.Range('D18').Value = '''1.001'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('D19').Value = '''0.00001009'
$ws.Range('E19').Value = '  -2.95%  '
$ws.Range('D20').Value = '''17.13'
$ws.Range('E20').Value = '  -2.83%  '
$ws.Range('E21').Value = '  -0.90%  '
$ws.Range('D22').Value = '28.697.10'
$ws.Range('E22').Value = '  -3.00%  '
$ws.Range('E23').Value = '  -5.12%  '
$ws.Range('D24').Value = '''11.34'
$ws.Range('E24').Value = '  -4.86%  '
$ws.Range('D25').Value = '''2.140'
$ws.Range('E25').Value = '  -1.97%  '
$ws.Range('D26').Value = '2.112.34'
$ws.Range('E26').Value = '  -2.16%  '
$ws.Range('D27').Value = '''153.26'
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('D28').Value = '''19.49'
$ws.Range('E28').Value = '  -2.58%  '
$ws.Range('D29').Value = '''6.123'
$ws.Range('E29').Value = '  -5.18%  '
$ws.Range('D30').Value = '''2.026'
$ws.Range('E30').Value = '  -3.32%  '
$ws.Range('D31').Value = '''117.62'
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('D32').Value = '''0.9779'
$ws.Range('E32').Value = '  -3.73%  '
$ws.Range('D33').Value = '''0.09426'
$ws.Range('E33').Value = '  -2.08%  '
$ws.Range('D34').Value = '''5.373'
$ws.Range('E34').Value = '  -4.34%  '
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('D36').Value = '''1.353'
$ws.Range('E36').Value = '  -1.88%  '
$ws.Range('D37').Value = '''0.06172'
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('D38').Value = '''0.02204'
$ws.Range('E38').Value = '  -3.79%  '
$ws.Range('D39').Value = '''1.165'
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.5729'
$ws.Range('E40').Value = '  -3.83%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '''7.609'
$ws.Range('E41').Value = '  -3.35%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '''10.17'
$ws.Range('E42').Value = '  -5.48%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '''0.1800'
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '''2.367'
$ws.Range('E44').Value = '  -3.53%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '''1.252'
$ws.Range('E45').Value = '  -2.77%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.5408'
$ws.Range('E46').Value = '  -2.85%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''11.85'
$ws.Range('E47').Value = '  -4.52%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.07155'
$ws.Range('E48').Value = '  -4.57%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '''1.916'
$ws.Range('E49').Value = '  -1.74%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '''114.22'
$ws.Range('E50').Value = '  -4.22%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = '''43.24'
$ws.Range('E51').Value = '  +3.01%  '
